$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the filename string in A2
$ws.Range("A2").Value = "1__211006132800_Waves_001.txt"

# Update the numeric values in row 2 (columns D through AK)
$newValues = @{
    "D2"  = 5
    "E2"  = 12
    "F2"  = 510.58
    "G2"  = 7.71
    "H2"  = 5.79
    "I2"  = 0.5
    "J2"  = 0.06
    "K2"  = 3.56
    "L2"  = 0.33
    "M2"  = 0.04
    "N2"  = 3.89
    "O2"  = 0.34
    "P2"  = 0.04
    "Q2"  = 42.83
    "R2"  = 7.14
    "S2"  = 0.91
    "T2"  = 2.79
    "U2"  = 0.37
    "V2"  = 0.05
    "W2"  = 131.24
    "X2"  = 12.59
    "Y2"  = 1.61
    "Z2"  = 8.41
    "AA2" = 0.8
    "AB2" = 0.1
    "AC2" = 7.44
    "AD2" = 0.64
    "AE2" = 0.08
    "AF2" = 11.1
    "AG2" = 0.46
    "AH2" = 0.06
    "AI2" = 17.34
    "AJ2" = 0.85
    "AK2" = 0.11
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}
